# Apply updated cryptocurrency market data (price & 1h volume change)
# to the cryptos worksheet, matching the scheduled GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '''25.982.29'
$ws.Range("E2").Value = '  -0.21%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '''1.629.09'
$ws.Range("E3").Value = '  -0.91%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.08%  '

# Row 5: BNB
$ws.Range("D5").Value = '''214.14'
$ws.Range("E5").Value = '  -0.96%  '

# Row 6: XRP
$ws.Range("E6").Value = '  -0.84%  '

# Row 7: USDC
$ws.Range("D7").Value = '''1.01'
$ws.Range("E7").Value = '  -0.07%  '

# Row 8: Cardano
$ws.Range("E8").Value = '  -1.89%  '

# Row 9: Dogecoin
$ws.Range("E9").Value = '  -3.13%  '

# Row 10: Solana
$ws.Range("D10").Value = '''18.46'
$ws.Range("E10").Value = '  -5.74%  '

# Row 11: TRON
$ws.Range("E11").Value = '  -1.30%  '

# Row 12: WrappedEther
$ws.Range("B12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D12").Value = '''1.857.03'
$ws.Range("E12").Value = '  -0.79%  '

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''4.19'
$ws.Range("E13").Value = '  -1.99%  '

# Row 14: Polkadot
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '''1.623.14'
$ws.Range("E14").Value = '  -2.06%  '

# Row 15: Polygon
$ws.Range("D15").Value = '''0.529'
$ws.Range("E15").Value = '  -2.91%  '

# Row 16: WrappedBTC
$ws.Range("D16").Value = '''25.974.20'
$ws.Range("E16").Value = '  -0.34%  '

# Row 17: ShibaInu
$ws.Range("D17").Value = '''0.0₃0740'
$ws.Range("E17").Value = '  -3.17%  '

# Row 18: Litecoin
$ws.Range("D18").Value = '''61.36'
$ws.Range("E18").Value = '  -3.31%  '

# Row 19: Dai
$ws.Range("E19").Value = '  -0.08%  '

# Row 20: BitcoinCash
$ws.Range("D20").Value = '''192.34'
$ws.Range("E20").Value = '  -1.16%  '

# Row 21: Uniswap
$ws.Range("D21").Value = '''4.25'
$ws.Range("E21").Value = '  -2.68%  '

# Row 22: Avalanche
$ws.Range("D22").Value = '''9.59'
$ws.Range("E22").Value = '  -3.51%  '

# Row 23: Chainlink
$ws.Range("E23").Value = '  -2.03%  '

# Row 24: Stellar
$ws.Range("E24").Value = '  +0.65%  '

# Row 25: Monero
$ws.Range("D25").Value = '''143.69'
$ws.Range("E25").Value = '  +0.29%  '

# Row 26: BinanceUSD
$ws.Range("D26").Value = '''1.01'
$ws.Range("E26").Value = '  -0.05%  '

# Row 27: Toncoin
$ws.Range("D27").Value = '''1.73'
$ws.Range("E27").Value = '  -3.29%  '

# Row 28: Cosmos
$ws.Range("D28").Value = '''6.74'
$ws.Range("E28").Value = '  -2.11%  '

# Row 29: EthereumClassic
$ws.Range("D29").Value = '''15.20'
$ws.Range("E29").Value = '  -2.07%  '

# Row 30: PancakeSwap
$ws.Range("E30").Value = '  -1.42%  '

# Row 31: Hedera
$ws.Range("E31").Value = '  -2.19%  '

# Row 32: Filecoin
$ws.Range("E32").Value = '  -4.05%  '

# Row 33: InternetComputer(DFINITY)
$ws.Range("E33").Value = '  -5.40%  '

# Row 34: LidoDAOToken
$ws.Range("E34").Value = '  -2.75%  '

# Row 35: HuobiToken
$ws.Range("E35").Value = '  -2.86%  '

# Row 36: Maker
$ws.Range("D36").Value = '''1.129.01'

# Row 37: ARBITRUM
$ws.Range("E37").Value = '  -5.67%  '

# Row 38: MXToken
$ws.Range("E38").Value = '  -1.45%  '

# Row 39: ImmutableX
$ws.Range("E39").Value = '  -3.55%  '

# Row 40: VeChain
$ws.Range("E40").Value = '  -2.41%  '

# Row 41: Quant
$ws.Range("D41").Value = '''98.16'
$ws.Range("E41").Value = '  -0.98%  '

# Row 42: RocketPoolETH
$ws.Range("D42").Value = '''1.765.79'
$ws.Range("E42").Value = '  -0.83%  '

# Row 43: TrustWalletToken
$ws.Range("D43").Value = '''0.764'
$ws.Range("E43").Value = '  -4.33%  '

# Row 44: FraxShare
$ws.Range("E44").Value = '  -5.42%  '

# Row 45: BabyDogeCoin
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").Value = '''0.0532'
$ws.Range("E45").Value = '  +1.95%  '

# Row 46: Cronos
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '''54.44'
$ws.Range("E46").Value = '  -3.71%  '

# Row 47: Aave
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '''1.48'
$ws.Range("E47").Value = '  -0.70%  '

# Row 48: RenderToken
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '''0.0₇0987'
$ws.Range("E48").Value = '  -15.93%  '

# Row 49: Mantle
$ws.Range("E49").Value = '  -0.41%  '

# Row 50: USDD
$ws.Range("D50").Value = '''1.01'
$ws.Range("E50").Value = '  +0.16%  '

# Row 51: EnergySwap
$ws.Range("E51").Value = '  -3.57%  '
